# Update "想去人数" (F column) counters and sold-out status (G column)
# across the 展览 / 演出 / 本地生活 / 全部类型 sheets, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2311
$ws1.Range("F7").Value = 497
$ws1.Range("F9").Value = 733
$ws1.Range("F11").Value = 725
$ws1.Range("F16").Value = 991
$ws1.Range("F17").Value = 17580
$ws1.Range("G17").Value = "已售罄"
$ws1.Range("F18").Value = 394
$ws1.Range("F19").Value = 42
$ws1.Range("F20").Value = 180
$ws1.Range("F25").Value = 156

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 193
$ws2.Range("F8").Value = 3318
$ws2.Range("F10").Value = 36
$ws2.Range("F16").Value = 2710

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 71
$ws3.Range("F4").Value = 514

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 71
$ws4.Range("F6").Value = 2311
$ws4.Range("F7").Value = 514
$ws4.Range("F12").Value = 497
$ws4.Range("F16").Value = 193
$ws4.Range("F18").Value = 733
$ws4.Range("F20").Value = 725
$ws4.Range("F25").Value = 991
$ws4.Range("F26").Value = 17580
$ws4.Range("G26").Value = "已售罄"
$ws4.Range("F28").Value = 3318
$ws4.Range("F30").Value = 36
$ws4.Range("F32").Value = 394
$ws4.Range("F33").Value = 42
$ws4.Range("F34").Value = 180
$ws4.Range("F43").Value = 156
$ws4.Range("F47").Value = 2710
